# Applies the StructureDefinition-see-also.xlsx content update:
#  - rebrand ibm.com -> linuxforhealth.org in canonical URL(s)
#  - bump Version 7.0.0 -> 8.0.0
#  - bump Date 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
#  - rebrand Publisher Alvearie Team -> LinuxForHealth Team
#  - clear the stale constraint text that used to sit on the top-level
#    Extension row (it now only applies to the Extension.extension row)
#  - rebrand the ibm.com extension-definition URLs used in the Elements
#    sheet's Type(s) column for the match/split extensions

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/see-also"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
# The top-level "Extension" row no longer carries the ele-1/ext-1
# constraint text in its "Constraint(s)" column (AI) - it moved to the
# Extension.extension row instead.
$wsElem.Range("AI2").Value = ""

# The "Fixed Value" for Extension.url mirrors the StructureDefinition's
# own canonical URL.
$wsElem.Range("Q12").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/see-also"

# Rebrand the extension-definition URLs referenced in the Type(s) column.
$wsElem.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-source-reference}
"
$wsElem.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-method}
"
$wsElem.Range("J7").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/split-method}
"
$wsElem.Range("J8").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-confidence-level}
"
$wsElem.Range("J9").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-confidence-score}
"
$wsElem.Range("J10").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-period}
"
$wsElem.Range("J11").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/match-detail}
"
